$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.930.55'
$ws.Range("E2").Value = '  -1.66%  '

$ws.Range("D3").Value = '2.347.54'
$ws.Range("E3").Value = '  -2.18%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.47%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0912'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.109'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.964'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.52%  '

$ws.Range("D15").Value = '2.703.57'
$ws.Range("E15").Value = '  -2.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.24%  '

$ws.Range("D17").Value = '2.375.88'
$ws.Range("E17").Value = '  -1.14%  '

$ws.Range("D18").Value = '44.954.38'
$ws.Range("E18").Value = '  -1.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +12.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000105'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '256.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.90%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.66%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0952'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.88'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.57'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '166.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.48%  '

$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.131'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.76%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.68'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0349'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.77%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.874.68'
$ws.Range("E42").Value = '  +13.44%  '

$ws.Range("B43").Value = 'BitcoinSV'
$ws.Range("C43").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '95.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.63%  '

$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '68.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.43%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.24%  '

$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.98%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.225'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.94%  '

$ws.Range("B48").Value = 'Celestia'
$ws.Range("C48").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.67%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.22%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.28%  '
